$wb = $excel.ActiveWorkbook

# --- Add the new "Spain" sheet by duplicating the existing "Italy" sheet ---
# (Spain's layout/styles/merges are identical to Italy's; only a couple of
# data cells and the sheet name differ.)
$italy = $wb.Worksheets.Item("Italy")
$italy.Copy($null, $italy)

$spain = $wb.Worksheets.Item($wb.Worksheets.Count)
$spain.Name = "Spain"

# --- Update the market-specific data cells on the new sheet ---
$spain.Range("B2").Value = "Spain Market"
$spain.Range("B4").Value = "NGC-3103/T2041"

# --- Restore Italy's view state (it is no longer the active tab) ---
$italy.Range("A1:D12").Select()

# --- Set view state / active cell on the new Spain sheet & activate it ---
$spain.Range("F15").Select()
$spain.Activate()
